$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 730, shifting the existing rows 730:806 down to 732:808
$ws.Range("A730:R731").EntireRow.Insert()

# New row 730
$ws.Range("A730").Value = 7
$ws.Range("B730").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C730").Value = "Ñuble"
$ws.Range("D730").Value = [DateTime]"2023-07-25"
$ws.Range("E730").Value = 16
$ws.Range("F730").Value = 100112020
$ws.Range("G730").Value = "Tomate"
$ws.Range("H730").Value = "Larga vida"
$ws.Range("I730").Value = "Primera"
$ws.Range("J730").Value = 100
$ws.Range("K730").Value = 19000
$ws.Range("L730").Value = 19000
$ws.Range("M730").Value = 19000
$ws.Range("N730").Value = "`$/bandeja 20 kilos"
$ws.Range("O730").Value = "Región de Arica y Parinacota"
$ws.Range("P730").Value = 950
$ws.Range("Q730").Value = 20
$ws.Range("R730").Value = "Hortaliza"

# New row 731
$ws.Range("A731").Value = 7
$ws.Range("B731").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C731").Value = "Ñuble"
$ws.Range("D731").Value = [DateTime]"2023-07-25"
$ws.Range("E731").Value = 16
$ws.Range("F731").Value = 100112020
$ws.Range("G731").Value = "Tomate"
$ws.Range("H731").Value = "Larga vida"
$ws.Range("I731").Value = "Segunda"
$ws.Range("J731").Value = 80
$ws.Range("K731").Value = 17000
$ws.Range("L731").Value = 17000
$ws.Range("M731").Value = 17000
$ws.Range("N731").Value = "`$/bandeja 20 kilos"
$ws.Range("O731").Value = "Región de Arica y Parinacota"
$ws.Range("P731").Value = 850
$ws.Range("Q731").Value = 20
$ws.Range("R731").Value = "Hortaliza"
